$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update existing task description (row 27, column A)
$ws.Range("A27").Value = "agregar en seguimiento de OT el estado del OT buscada y posibilidad de imprimir"

# Add new row 28
$ws.Range("A28").Value = "implementar en tabla cobranza el importe total de la cobranza y terminar reporte de saldo deudor"
$ws.Range("B28").Value = "en proceso"

# Update view to reflect scrolled/selected state
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("C29").Select()
